$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.316781520843506
$ws.Range("B1").Value = 1.865177392959595
$ws.Range("C1").Value = 1.790903806686401
$ws.Range("D1").Value = 4.951620101928711
$ws.Range("E1").Value = 1.321988105773926
